$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 57: becomes the old row 60 record (Knärot, count 17), with coordinates rounded to integers
$ws.Range("A57").Value = 112145544
$ws.Range("I57").Value = "17"
$ws.Range("Q57").Value = 653024
$ws.Range("R57").Value = 6675364
$ws.Range("Z57").Value = "11:03"
$ws.Range("AB57").Value = "11:04"

# Row 58: unchanged record, only coordinates get rounded to integers
$ws.Range("Q58").Value = 653012
$ws.Range("R58").Value = 6675152

# Row 59: becomes the old row 57 record (Knärot, count 14), with coordinates rounded to integers
$ws.Range("A59").Value = 112145545
$ws.Range("B59").Value = 96348
$ws.Range("D59").Value = "VU"
$ws.Range("E59").Value = 220787
$ws.Range("F59").Value = "Knärot"
$ws.Range("G59").Value = "Goodyera repens"
$ws.Range("H59").Value = "(L.) R. Br."
$ws.Range("I59").Value = "14"
$ws.Range("J59").Value = "plantor/tuvor"
$ws.Range("Q59").Value = 653038
$ws.Range("R59").Value = 6675341
$ws.Range("Z59").Value = "11:07"
$ws.Range("AB59").Value = "11:08"

# Row 60: becomes the old row 59 record (Brandticka), with coordinates rounded to integers
$ws.Range("A60").Value = 112145539
$ws.Range("B60").Value = 90018
$ws.Range("D60").Value = "LC"
$ws.Range("E60").Value = 1339
$ws.Range("F60").Value = "Brandticka"
$ws.Range("G60").Value = "Pycnoporellus fulgens"
$ws.Range("H60").Value = "(Fr.) Donk"
$ws.Range("I60").Value = ""
$ws.Range("J60").Value = ""
$ws.Range("Q60").Value = 652997
$ws.Range("R60").Value = 6675310
$ws.Range("Z60").Value = "10:44"
$ws.Range("AB60").Value = "10:44"
